$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the Romansh ("rm") language to the "classes" sheet by inserting two new
# columns:
#  - a new "rm" language-code column before the old column F
#  - a new "comment_rm" column before the old column J (which becomes L)
$ws.Columns.Item(6).Insert() | Out-Null
$ws.Columns.Item(11).Insert() | Out-Null

# Header row
$ws.Range("F1").Value = "rm"
$ws.Range("K1").Value = "comment_rm"

# Data rows 2-7: Romansh language name
$ws.Range("F2:F7").Value = "Rumantsch"
$ws.Range("K2:K7").Value = "Rumantsch"

# A few rows in this table are shaded (style carried by columns B/C); make the
# new "rm" cells in those rows match that shading.
$ws.Range("B4").Copy() | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null
$ws.Range("B5").Copy() | Out-Null
$ws.Range("F5").PasteSpecial(-4122) | Out-Null
$ws.Range("B7").Copy() | Out-Null
$ws.Range("F7").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false | Out-Null

# Move the active tab/selection from "Owner" to "classes" (this sheet),
# matching the saved workbook view state.
$ws.Range("K13").Select() | Out-Null
$ws.Activate() | Out-Null
